# Applies the numeric restatements captured in the scheduled-runner diff
# for Sheets/Bahamut_Profits.xlsx across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR tabs.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1778.5122
$ws.Range("I132").Value = 2066.7932
$ws.Range("J132").Value = 1081.8334
$ws.Range("K132").Value = 6200.3796
$ws.Range("L132").Value = 3245.5002
$ws.Range("M132").Value = -3670.3796
$ws.Range("N132").Value = -8305.5002

$ws.Range("H137").Value = 770.7843
$ws.Range("I137").Value = 698.6896400000001
$ws.Range("J137").Value = 865.8182
$ws.Range("K137").Value = 2096.06892
$ws.Range("L137").Value = 2597.4546
$ws.Range("M137").Value = 453.9310799999998
$ws.Range("N137").Value = -7697.4546

$ws.Range("H141").Value = 1923.585
$ws.Range("I141").Value = 687.2727
$ws.Range("J141").Value = 7967.778
$ws.Range("K141").Value = 2061.8181
$ws.Range("L141").Value = 23903.334
$ws.Range("M141").Value = 3118.1819
$ws.Range("N141").Value = -34263.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7460.6333
$ws.Range("I32").Value = 7578.672
$ws.Range("J32").Value = 7170.077
$ws.Range("K32").Value = 7578.672
$ws.Range("L32").Value = 7170.077
$ws.Range("M32").Value = -7291.672
$ws.Range("N32").Value = -7744.077

$ws.Range("H61").Value = 1049.8108
$ws.Range("I61").Value = 933.8570999999999
$ws.Range("J61").Value = 1202
$ws.Range("K61").Value = 933.8570999999999
$ws.Range("L61").Value = 1202
$ws.Range("M61").Value = -721.8570999999999
$ws.Range("N61").Value = -1626

$ws.Range("H132").Value = 1334.4073
$ws.Range("I132").Value = 1188.125
$ws.Range("J132").Value = 2504.6667
$ws.Range("K132").Value = 3564.375
$ws.Range("L132").Value = 7514.000100000001
$ws.Range("M132").Value = -1034.375
$ws.Range("N132").Value = -12574.0001

$ws.Range("H136").Value = 1049.8108
$ws.Range("I136").Value = 933.8570999999999
$ws.Range("J136").Value = 1202
$ws.Range("K136").Value = 2801.5713
$ws.Range("L136").Value = 3606
$ws.Range("M136").Value = -251.5712999999996
$ws.Range("N136").Value = -8706

$ws.Range("H138").Value = 71151.14
$ws.Range("J138").Value = 71151.14
$ws.Range("L138").Value = 71151.14
$ws.Range("N138").Value = -81431.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 18803.69
$ws.Range("I134").Value = 1335.8043
$ws.Range("J134").Value = 85763.914
$ws.Range("K134").Value = 4007.4129
$ws.Range("L134").Value = 257291.742
$ws.Range("M134").Value = -1472.4129
$ws.Range("N134").Value = -262361.742

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2765.0981
$ws.Range("I31").Value = 2566.6572
$ws.Range("J31").Value = 3199.1875
$ws.Range("K31").Value = 2566.6572
$ws.Range("L31").Value = 3199.1875
$ws.Range("M31").Value = -2271.6572
$ws.Range("N31").Value = -3789.1875

$ws.Range("H34").Value = 2765.0981
$ws.Range("I34").Value = 2566.6572
$ws.Range("J34").Value = 3199.1875
$ws.Range("K34").Value = 2566.6572
$ws.Range("L34").Value = 3199.1875
$ws.Range("M34").Value = -2364.6572
$ws.Range("N34").Value = -3603.1875

$ws.Range("H58").Value = 1446.826
$ws.Range("I58").Value = 1815.2307
$ws.Range("J58").Value = 967.9
$ws.Range("K58").Value = 1815.2307
$ws.Range("L58").Value = 967.9
$ws.Range("M58").Value = -1612.2307
$ws.Range("N58").Value = -1373.9

$ws.Range("H108").Value = 30433.334
$ws.Range("I108").Value = 29850
$ws.Range("J108").Value = 30725
$ws.Range("K108").Value = 29850
$ws.Range("L108").Value = 30725
$ws.Range("M108").Value = -26010
$ws.Range("N108").Value = -38405

$ws.Range("H132").Value = 1707.1666
$ws.Range("I132").Value = 1272.5
$ws.Range("J132").Value = 2250.5
$ws.Range("K132").Value = 3817.5
$ws.Range("L132").Value = 6751.5
$ws.Range("M132").Value = -1287.5
$ws.Range("N132").Value = -11811.5

$ws.Range("H134").Value = 1396.836
$ws.Range("I134").Value = 1344.7551
$ws.Range("J134").Value = 1609.5
$ws.Range("K134").Value = 4034.2653
$ws.Range("L134").Value = 4828.5
$ws.Range("M134").Value = -1499.2653
$ws.Range("N134").Value = -9898.5

$ws.Range("H136").Value = 1446.826
$ws.Range("I136").Value = 1815.2307
$ws.Range("J136").Value = 967.9
$ws.Range("K136").Value = 5445.6921
$ws.Range("L136").Value = 2903.7
$ws.Range("M136").Value = -2895.6921
$ws.Range("N136").Value = -8003.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 500.45
$ws.Range("I122").Value = 235.9375
$ws.Range("J122").Value = 676.7917
$ws.Range("K122").Value = 2123.4375
$ws.Range("L122").Value = 6091.1253
$ws.Range("M122").Value = 326.5625
$ws.Range("N122").Value = -10991.1253

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 5159.6665
$ws.Range("I53").Value = 4995.6
$ws.Range("K53").Value = 4995.6
$ws.Range("M53").Value = -4364.6

$ws.Range("H102").Value = 1302.1818
$ws.Range("I102").Value = 1305.7778
$ws.Range("J102").Value = 1286
$ws.Range("K102").Value = 1305.7778
$ws.Range("L102").Value = 1286
$ws.Range("M102").Value = 316.2221999999999
$ws.Range("N102").Value = -4530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1583.4445
$ws.Range("I16").Value = 1443.5834
$ws.Range("J16").Value = 1863.1666
$ws.Range("K16").Value = 1443.5834
$ws.Range("L16").Value = 1863.1666
$ws.Range("M16").Value = -1273.5834
$ws.Range("N16").Value = -2203.1666

$ws.Range("H40").Value = 1123878.9
$ws.Range("I40").Value = 2021382
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 2021382
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -2021246
$ws.Range("N40").Value = -2272

$ws.Range("H46").Value = 2485
$ws.Range("I46").Value = 2156
$ws.Range("J46").Value = 2634.5454
$ws.Range("K46").Value = 2156
$ws.Range("L46").Value = 2634.5454
$ws.Range("M46").Value = -1968
$ws.Range("N46").Value = -3010.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 733.8333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 733.8333
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2201.4999
$ws.Range("N113").Value = -6541.4999
$ws.Range("M113").ClearContents()

$ws.Range("H132").Value = 816.3200000000001
$ws.Range("I132").Value = 803.6
$ws.Range("J132").Value = 867.2
$ws.Range("K132").Value = 2410.8
$ws.Range("L132").Value = 2601.6
$ws.Range("M132").Value = 119.1999999999998
$ws.Range("N132").Value = -7661.6

$ws.Range("H136").Value = 1704.1666
$ws.Range("I136").Value = 2133.75
$ws.Range("J136").Value = 845
$ws.Range("K136").Value = 6401.25
$ws.Range("L136").Value = 2535
$ws.Range("M136").Value = -3851.25
$ws.Range("N136").Value = -7635

